# Auto-generated edit script applying numeric updates to Goblin_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 5000
$ws.Range("I10").Value = 5000
$ws.Range("K10").Value = 5000
$ws.Range("M10").Value = -4707
$ws.Range("H40").Value = 3556.5715
$ws.Range("I40").Value = 1800.25
$ws.Range("J40").Value = 3969.8235
$ws.Range("K40").Value = 1800.25
$ws.Range("L40").Value = 3969.8235
$ws.Range("M40").Value = -1625.25
$ws.Range("N40").Value = -4319.8235
$ws.Range("H64").Value = 8424.678
$ws.Range("I64").Value = 4944.6665
$ws.Range("K64").Value = 4944.6665
$ws.Range("M64").Value = -4696.6665
$ws.Range("H67").Value = 8424.678
$ws.Range("I67").Value = 4944.6665
$ws.Range("K67").Value = 4944.6665
$ws.Range("M67").Value = -4086.6665
$ws.Range("H105").Value = 60399.6
$ws.Range("J105").Value = 60399.6
$ws.Range("L105").Value = 60399.6
$ws.Range("N105").Value = -67387.60000000001
$ws.Range("H137").Value = 4586.5864
$ws.Range("I137").Value = 5769.3335
$ws.Range("K137").Value = 17308.0005
$ws.Range("M137").Value = -14758.0005

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2851032.5
$ws.Range("I122").Value = 3586007.8
$ws.Range("J122").Value = 3003.375
$ws.Range("K122").Value = 10758023.4
$ws.Range("L122").Value = 9010.125
$ws.Range("M122").Value = -10755573.4
$ws.Range("N122").Value = -13910.125
$ws.Range("H123").Value = 90000
$ws.Range("J123").Value = 90000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -99800
$ws.Range("H132").Value = 1944.7291
$ws.Range("I132").Value = 1469.5641
$ws.Range("K132").Value = 4408.692300000001
$ws.Range("M132").Value = -1878.692300000001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3238.9
$ws.Range("I94").Value = 2869.2942
$ws.Range("J94").Value = 5333.3335
$ws.Range("K94").Value = 2869.2942
$ws.Range("L94").Value = 5333.3335
$ws.Range("M94").Value = -2418.2942
$ws.Range("N94").Value = -6235.3335
$ws.Range("H102").Value = 13708.333
$ws.Range("I102").Value = 5450
$ws.Range("J102").Value = 55000
$ws.Range("K102").Value = 5450
$ws.Range("L102").Value = 55000
$ws.Range("M102").Value = -2205
$ws.Range("N102").Value = -61490
$ws.Range("H132").Value = 53463.555
$ws.Range("J132").Value = 53463.555
$ws.Range("L132").Value = 53463.555
$ws.Range("N132").Value = -63583.555
$ws.Range("H134").Value = 2239.138
$ws.Range("I134").Value = 1473.238
$ws.Range("J134").Value = 4249.625
$ws.Range("K134").Value = 4419.714
$ws.Range("L134").Value = 12748.875
$ws.Range("M134").Value = -1884.714
$ws.Range("N134").Value = -17818.875
$ws.Range("H135").Value = 71999.60000000001
$ws.Range("J135").Value = 71999.60000000001
$ws.Range("L135").Value = 71999.60000000001
$ws.Range("N135").Value = -82139.60000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2677.4905
$ws.Range("I31").Value = 1081.7273
$ws.Range("K31").Value = 1081.7273
$ws.Range("M31").Value = -786.7273
$ws.Range("H34").Value = 2677.4905
$ws.Range("I34").Value = 1081.7273
$ws.Range("K34").Value = 1081.7273
$ws.Range("M34").Value = -879.7273
$ws.Range("H132").Value = 3815
$ws.Range("J132").Value = 2506.75
$ws.Range("L132").Value = 7520.25
$ws.Range("N132").Value = -12580.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1353028.1
$ws.Range("I4").Value = 669070.7
$ws.Range("K4").Value = 2007212.1
$ws.Range("M4").Value = -2007100.1
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H49").Value = 1100.5
$ws.Range("J49").Value = 999.5
$ws.Range("L49").Value = 2998.5
$ws.Range("N49").Value = -3310.5
$ws.Range("H114").Value = 743
$ws.Range("I114").Value = 949
$ws.Range("J114").Value = 331
$ws.Range("K114").Value = 2847
$ws.Range("L114").Value = 993
$ws.Range("M114").Value = 407
$ws.Range("N114").Value = -7501
$ws.Range("H117").Value = 1755.75
$ws.Range("I117").Value = 2317.5
$ws.Range("J117").Value = 1194
$ws.Range("K117").Value = 6952.5
$ws.Range("L117").Value = 3582
$ws.Range("M117").Value = -3510.5
$ws.Range("N117").Value = -10466
$ws.Range("H131").Value = 2300924
$ws.Range("I131").Value = 857.5
$ws.Range("K131").Value = 2572.5
$ws.Range("M131").Value = 2467.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4085.3845
$ws.Range("I132").Value = 4085.3845
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12256.1535
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9726.1535
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 71808
$ws.Range("J141").Value = 71808
$ws.Range("L141").Value = 71808
$ws.Range("N141").Value = -82168

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 9975
$ws.Range("I13").Value = 6950
$ws.Range("J13").Value = 13000
$ws.Range("K13").Value = 6950
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = -6810
$ws.Range("N13").Value = -13280
$ws.Range("H46").Value = 854.5833
$ws.Range("I46").Value = 781.375
$ws.Range("K46").Value = 781.375
$ws.Range("M46").Value = -593.375
$ws.Range("H74").Value = 46666.332
$ws.Range("I74").Value = 39999.5
$ws.Range("J74").Value = 60000
$ws.Range("K74").Value = 39999.5
$ws.Range("L74").Value = 60000
$ws.Range("M74").Value = -39001.5
$ws.Range("N74").Value = -61996
$ws.Range("H77").Value = 46666.332
$ws.Range("I77").Value = 39999.5
$ws.Range("J77").Value = 60000
$ws.Range("K77").Value = 119998.5
$ws.Range("L77").Value = 180000
$ws.Range("M77").Value = -115006.5
$ws.Range("N77").Value = -189984
$ws.Range("H122").Value = 4935.1665
$ws.Range("J122").Value = 6280.8
$ws.Range("L122").Value = 18842.4
$ws.Range("N122").Value = -23742.4

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 1500
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766
$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652
$ws.Range("H114").Value = 87398
$ws.Range("J114").Value = 87398
$ws.Range("L114").Value = 87398
$ws.Range("N114").Value = -96076
$ws.Range("H122").Value = 4584.3
$ws.Range("J122").Value = 7200
$ws.Range("L122").Value = 21600
$ws.Range("N122").Value = -26500
$ws.Range("H132").Value = 6655.6787
$ws.Range("I132").Value = 6778.1816
$ws.Range("J132").Value = 6576.4116
$ws.Range("K132").Value = 20334.5448
$ws.Range("L132").Value = 19729.2348
$ws.Range("M132").Value = -17804.5448
$ws.Range("N132").Value = -24789.2348
